# Update the "version" value on the TestData sheet from 107.0 to 112.0
# and move the active selection to G13, matching the committed change.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TestData")

# Leading apostrophe keeps this a text value (matches the original cell,
# which is stored as a shared string with a quote-prefix style) instead of
# Excel auto-converting "112.0" into the number 112.
$ws.Range("D2").Value = "'112.0"

$ws.Activate()
$ws.Range("G13").Select()
